$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND (replace): $old"
    }
}

function Italicize-Suffix($contextPhrase, $word) {
    # Finds a unique, longer phrase that ends with $word, then italicizes
    # just the trailing $word (by character offset), leaving the rest of
    # the phrase untouched.
    $r = $d.Content
    $ok = $r.Find.Execute($contextPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND (italicize context): $contextPhrase"
        return
    }
    $wordLen = $word.Length
    $start = $r.End - $wordLen
    $end = $r.End
    $sub = $d.Range($start, $end)
    if ($sub.Text -ne $word) {
        Write-Output "MISMATCH (italicize): expected [$word] got [$($sub.Text)]"
        return
    }
    $sub.Font.Italic = 1
}

# --- Hunk 1: cross-sector interactions paragraph ---
Replace-Text "other sectors. Several different kinds of cross-sector interactions can be defined as having important impact" "stakeholders in other sectors. Several different kinds of cross-sector interactions can be defined as having an important impact"

# --- Hunk 2: James Holdren date ---
Replace-Text "James Holdren (in 1) and" "James Holdren (in 2013) and"

# --- Hunk 3: extend "how standards arise" sentence ---
Replace-Text "how standards arise so that these goals are achieved." "how standards arise so that these goals are achieved. Importantly, open-source standards seem to well-match at least some of these characteristics."

# --- Hunk 4a: prepend new sentence before "A compelling road map" ---
Replace-Text "A compelling road map towards implementation and adoption of" "The other side of policies is the implementation of these policies in practice by developers of open-source standards and by the communities to which the standards pertain. A compelling road map towards implementation and adoption of"

# --- Hunk 4b: community-developed standards -> open science practices ... ---
Replace-Text "community-developed standards is offered in a blog post authored by the Center for Open Science’s Brian Nosek, entitled" "open science practices in general and open-source standards in particular is offered in a blog post authored by the Center for Open Science’s co-founder and executive director, Brian Nosek, entitled"

# --- Hunk 4c: rewrite the final sentence(s) of the paragraph ---
Replace-Text "these pieces, which make adoption of standards possible, and maybe even easy, and the policy goals, arises from a community of practice that makes adoption of standards normative. Once all of these pieces are in place, making adoption of open science standards required becomes more straightforward and less onerous." "these pieces, which makes the adoption of standards possible, and maybe even easy, and the policy goals, arises from a community of practice that makes adoption of standards normative. Once all of these pieces are in place, making adoption of open science standards required through policy becomes more straightforward and less onerous."

# --- Hunk 4d: italicize "normative" and "required" (the new occurrences only) ---
Italicize-Suffix "adoption of standards normative" "normative"
Italicize-Suffix "open science standards required" "required"

# --- Hunk 5: extend BIDS sentence ---
Replace-Text "the Brain Imaging Data Structure standard in neuroscience." "the Brain Imaging Data Structure standard in neuroscience. Where large governmental funding agencies may not have"

Write-Output "DONE"
